$d = $word.ActiveDocument

# 1) Remove the anchored picture ("Picture 1") that was wrapped in the
#    first run of the title paragraph.
if ($d.Shapes.Count -gt 0) {
    $shp = $d.Shapes.Item(1)
    $shp.Delete()
}

# 2) Split the title paragraph so a new, empty, centered paragraph
#    (inheriting the same sz/szCs run formatting) follows it. Using
#    Find/Replace with the special "^p" sequence splits the paragraph
#    in place without synthesizing a stray empty run in the new
#    paragraph (unlike Range.InsertParagraphAfter / TypeParagraph).
$d.Content.Find.Execute("Power Failure Management", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Power Failure Management^p", 2)
